$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 242.66667
$ws.Range("I38").Value = 242.66667
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 728.00001
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -356.00001
$ws.Range("N38").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 50864.812
$ws.Range("I98").Value = 2151.25
$ws.Range("J98").Value = 99578.375
$ws.Range("K98").Value = 2151.25
$ws.Range("L98").Value = 99578.375
$ws.Range("M98").Value = -653.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 50864.812
$ws.Range("I122").Value = 2151.25
$ws.Range("J122").Value = 99578.375
$ws.Range("K122").Value = 6453.75
$ws.Range("L122").Value = 298735.125
$ws.Range("M122").Value = -4003.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1674
$ws.Range("I125").Value = 1865.3334
$ws.Range("J125").Value = 1559.2
$ws.Range("K125").Value = 16788.0006
$ws.Range("L125").Value = 14032.8
$ws.Range("M125").Value = -14328.0006
$ws.Range("N125").Value = -18952.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 16590.469
$ws.Range("I129").Value = 1028.3334
$ws.Range("J129").Value = 21354.389
$ws.Range("K129").Value = 3085.0002
$ws.Range("L129").Value = 64063.167
$ws.Range("M129").Value = 1914.9998
$ws.Range("N129").Value = -74063.167

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 31771.203
$ws.Range("I132").Value = 22819.857
$ws.Range("J132").Value = 61012.266
$ws.Range("K132").Value = 68459.571
$ws.Range("L132").Value = 183036.798
$ws.Range("M132").Value = -65929.571
$ws.Range("N132").Value = -188096.798

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1331581
$ws.Range("I137").Value = 2028755.6
$ws.Range("J137").Value = 6949.2
$ws.Range("K137").Value = 6086266.800000001
$ws.Range("L137").Value = 20847.6
$ws.Range("M137").Value = -6083716.800000001
$ws.Range("N137").Value = -25947.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 464.44446
$ws.Range("I4").Value = 463.33334
$ws.Range("J4").Value = 466.66666
$ws.Range("K4").Value = 463.33334
$ws.Range("L4").Value = 466.66666
$ws.Range("M4").Value = -347.33334
$ws.Range("N4").Value = -698.66666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 38626.5
$ws.Range("I23").Value = 70006
$ws.Range("J23").Value = 28166.666
$ws.Range("K23").Value = 70006
$ws.Range("L23").Value = 28166.666
$ws.Range("M23").Value = -69747
$ws.Range("N23").Value = -28684.666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11500.635
$ws.Range("I32").Value = 10262.915
$ws.Range("J32").Value = 29757
$ws.Range("K32").Value = 10262.915
$ws.Range("L32").Value = 29757
$ws.Range("M32").Value = -9975.915000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 37703.332
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 37703.332
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 37703.332
$ws.Range("N37").Value = -38249.332
$ws.Range("M37").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 37813.332
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 37813.332
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 37813.332
$ws.Range("N44").Value = -38789.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1895.7858
$ws.Range("I122").Value = 1673.2858
$ws.Range("J122").Value = 2118.2856
$ws.Range("K122").Value = 5019.857400000001
$ws.Range("L122").Value = 6354.8568
$ws.Range("M122").Value = -2569.857400000001
$ws.Range("N122").Value = -11254.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 10418153
$ws.Range("I132").Value = 13890001
$ws.Range("J132").Value = 2607
$ws.Range("K132").Value = 41670003
$ws.Range("L132").Value = 7821
$ws.Range("M132").Value = -41667473
$ws.Range("N132").Value = -12881

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 20000
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 20000
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 20000
$ws.Range("N9").Value = -20336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 10000
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 10000
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 10000
$ws.Range("N45").Value = -11186

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 740431.5600000001
$ws.Range("I132").Value = 1983.9231
$ws.Range("J132").Value = 2340401.5
$ws.Range("K132").Value = 5951.7693
$ws.Range("L132").Value = 7021204.5
$ws.Range("M132").Value = -3421.7693

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5940.76
$ws.Range("I3").Value = 2144.4443
$ws.Range("J3").Value = 8076.1875
$ws.Range("K3").Value = 6433.3329
$ws.Range("L3").Value = 24228.5625
$ws.Range("M3").Value = -6321.3329
$ws.Range("N3").Value = -24452.5625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 273.35
$ws.Range("I12").Value = 148.63158
$ws.Range("J12").Value = 386.1905
$ws.Range("K12").Value = 445.8947400000001
$ws.Range("L12").Value = 1158.5715
$ws.Range("M12").Value = -272.8947400000001
$ws.Range("N12").Value = -1504.5715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 37921.25
$ws.Range("I18").Value = 43295.715
$ws.Range("J18").Value = 300
$ws.Range("K18").Value = 129887.145
$ws.Range("L18").Value = 900
$ws.Range("M18").Value = -129718.145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 666.6667
$ws.Range("I32").Value = 600
$ws.Range("J32").Value = 700
$ws.Range("K32").Value = 1800
$ws.Range("L32").Value = 2100
$ws.Range("M32").Value = -1517
$ws.Range("N32").Value = -2666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1246.7433
$ws.Range("I68").Value = 1051.6666
$ws.Range("J68").Value = 1309.4464
$ws.Range("K68").Value = 3154.9998
$ws.Range("L68").Value = 3928.3392
$ws.Range("M68").Value = -2343.9998
$ws.Range("N68").Value = -5550.3392

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1246.7433
$ws.Range("I71").Value = 1051.6666
$ws.Range("J71").Value = 1309.4464
$ws.Range("K71").Value = 9464.999400000001
$ws.Range("L71").Value = 11785.0176
$ws.Range("M71").Value = -5408.999400000001
$ws.Range("N71").Value = -19897.0176

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 681.0700000000001
$ws.Range("I107").Value = 324.28
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 972.8399999999999
$ws.Range("L107").Value = 2400
$ws.Range("M107").Value = 947.1600000000001
$ws.Range("N107").Value = -6240

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2269.6924
$ws.Range("I122").Value = 1921.2
$ws.Range("J122").Value = 2487.5
$ws.Range("K122").Value = 5763.6
$ws.Range("L122").Value = 7462.5
$ws.Range("M122").Value = -3313.6
$ws.Range("N122").Value = -12362.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 11369.167
$ws.Range("I126").Value = 17179.428
$ws.Range("J126").Value = 3234.8
$ws.Range("K126").Value = 51538.284
$ws.Range("L126").Value = 9704.400000000001
$ws.Range("M126").Value = -49068.284
$ws.Range("N126").Value = -14644.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 501002
$ws.Range("I122").Value = 501002
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1503006
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1500556
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 20001592

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 29500
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 29500
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 29500
$ws.Range("N63").Value = -30748

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H66").Value = 29500
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 29500
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 88500
$ws.Range("N66").Value = -94740

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5002
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 5002
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 10004
$ws.Range("N81").Value = -12126
$ws.Range("M81").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 5002
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 5002
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 50020
$ws.Range("N84").Value = -60628
$ws.Range("M84").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 7144008
$ws.Range("I122").Value = 28571428
$ws.Range("J122").Value = 1534.3334
$ws.Range("K122").Value = 85714284
$ws.Range("L122").Value = 4603.0002
$ws.Range("M122").Value = -85711834
$ws.Range("N122").Value = -9503.0002
